$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.352.42"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.282.26"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.24"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.14"
$ws.Range("E6").Value = "  -1.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.648"
$ws.Range("E7").Value = "  +3.77%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.68"
$ws.Range("E10").Value = "  -2.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.21"
$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("E14").Value = "  -2.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.625.68"
$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.863"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.277.67"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.338.83"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.74"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.29"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.42"
$ws.Range("E22").Value = "  -1.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.58"
$ws.Range("E23").Value = "  +0.93%  "

$ws.Range("E24").Value = "  +4.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.45"
$ws.Range("E25").Value = "  -4.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.55"
$ws.Range("E27").Value = "  -1.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.96"
$ws.Range("E28").Value = "  -1.55%  "

$ws.Range("E30").Value = "  -1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.37"
$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0902"
$ws.Range("E33").Value = "  -2.83%  "

$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("E35").Value = "  +3.28%  "

$ws.Range("E36").Value = "  +3.59%  "

$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.84"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("E39").Value = "  -3.30%  "

$ws.Range("E40").Value = "  +8.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.27"
$ws.Range("E41").Value = "  +4.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.71"
$ws.Range("E42").Value = "  +2.86%  "

$ws.Range("E43").Value = "  -2.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.04"
$ws.Range("E44").Value = "  -4.15%  "

$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("E46").Value = "  -1.29%  "

$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.27"
$ws.Range("E47").Value = "  +3.68%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.57"
$ws.Range("E48").Value = "  -2.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  +0.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "100.30"
$ws.Range("E50").Value = "  -2.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.98"
$ws.Range("E51").Value = "  +29.06%  "
